# Update countries & provincias Spain
# - Reorder "Japon" to just after "Malasia" (before "Ecuador"), shifting the
#   rows that used to hold Ecuador/Filipinas/Pakistan down by one and giving
#   "Japon" fresh, updated figures.
# - Refresh the "Estados Unidos" (row 4) and "Brasil" (row 19) statistics.
# - Bump the "Datos actualizados..." timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 6 de Abril de 2020 a las 02:22"

# --- Estados Unidos (row 4) updated totals --------------------------------
$ws.Cells.Item(4, 2).Value2 = 336550   # Casos totales
$ws.Cells.Item(4, 3).Value2 = 25193    # Nuevos casos
$ws.Cells.Item(4, 5).Value2 = 309695   # Recuperados
$ws.Cells.Item(4, 7).Value2 = 1159     # Muertes hoy
$ws.Cells.Item(4, 8).Value2 = 9610     # Muertes

# --- Brasil (row 19) updated totals ---------------------------------------
$ws.Cells.Item(19, 2).Value2 = 11254   # Casos totales
$ws.Cells.Item(19, 3).Value2 = 894     # Nuevos casos
$ws.Cells.Item(19, 5).Value2 = 10641   # Recuperados

# --- Move "Japon" up above "Ecuador" (rows 34-37) --------------------------
# Row 34 becomes Japon with fresh numbers; former occupants of rows 34-36
# (Ecuador, Filipinas, Pakistan) shift down one row, keeping their own data.

# Row 34: Japon (new data)
$ws.Cells.Item(34, 1).Value2 = "Japon"
$ws.Cells.Item(34, 2).Value2 = 3654
$ws.Cells.Item(34, 3).Value2 = 515
$ws.Cells.Item(34, 4).Value2 = 575
$ws.Cells.Item(34, 5).Value2 = 2994
$ws.Cells.Item(34, 6).Value2 = 69
$ws.Cells.Item(34, 7).Value2 = 8
$ws.Cells.Item(34, 8).Value2 = 85

# Row 35: Ecuador (previous row-34 data)
$ws.Cells.Item(35, 1).Value2 = "Ecuador"
$ws.Cells.Item(35, 2).Value2 = 3646
$ws.Cells.Item(35, 3).Value2 = 181
$ws.Cells.Item(35, 4).Value2 = 100
$ws.Cells.Item(35, 5).Value2 = 3366
$ws.Cells.Item(35, 6).Value2 = 100
$ws.Cells.Item(35, 7).Value2 = 8
$ws.Cells.Item(35, 8).Value2 = 180

# Row 36: Filipinas (previous row-35 data)
$ws.Cells.Item(36, 1).Value2 = "Filipinas"
$ws.Cells.Item(36, 2).Value2 = 3246
$ws.Cells.Item(36, 3).Value2 = 152
$ws.Cells.Item(36, 4).Value2 = 64
$ws.Cells.Item(36, 5).Value2 = 3030
$ws.Cells.Item(36, 6).Value2 = 1
$ws.Cells.Item(36, 7).Value2 = 8
$ws.Cells.Item(36, 8).Value2 = 152

# Row 37: Pakistan (previous row-36 data)
$ws.Cells.Item(37, 1).Value2 = "Pakistan"
$ws.Cells.Item(37, 2).Value2 = 3157
$ws.Cells.Item(37, 3).Value2 = 339
$ws.Cells.Item(37, 4).Value2 = 211
$ws.Cells.Item(37, 5).Value2 = 2899
$ws.Cells.Item(37, 6).Value2 = 18
$ws.Cells.Item(37, 7).Value2 = 6
$ws.Cells.Item(37, 8).Value2 = 47
